$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two row swaps)
# Each cell is forced to Text storage (NumberFormat "@") to match the
# source data which stores prices/percentages as literal strings, then
# the style is reset to "Normal" so no stray style index is left behind.

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '69.915.29'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  +0.33%  '
$ws.Cells.Item(2, 5).Style = 'Normal'
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.690.88'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  -0.42%  '
$ws.Cells.Item(3, 5).Style = 'Normal'
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$ws.Cells.Item(4, 5).Style = 'Normal'
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '653.17'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  -3.57%  '
$ws.Cells.Item(5, 5).Style = 'Normal'
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '162.19'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  +0.24%  '
$ws.Cells.Item(6, 5).Style = 'Normal'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  +0.04%  '
$ws.Cells.Item(7, 5).Style = 'Normal'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  +0.51%  '
$ws.Cells.Item(8, 5).Style = 'Normal'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  -1.39%  '
$ws.Cells.Item(9, 5).Style = 'Normal'
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '7.18'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  +0.56%  '
$ws.Cells.Item(10, 5).Style = 'Normal'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  +0.35%  '
$ws.Cells.Item(11, 5).Style = 'Normal'
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.0000233'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '  -0.84%  '
$ws.Cells.Item(12, 5).Style = 'Normal'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '4.310.26'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '  -0.46%  '
$ws.Cells.Item(13, 5).Style = 'Normal'
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '32.85'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  +0.10%  '
$ws.Cells.Item(14, 5).Style = 'Normal'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '3.696.45'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  +0.18%  '
$ws.Cells.Item(15, 5).Style = 'Normal'
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '69.889.55'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  +0.24%  '
$ws.Cells.Item(16, 5).Style = 'Normal'
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  +0.36%  '
$ws.Cells.Item(17, 5).Style = 'Normal'
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '6.56'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  +0.86%  '
$ws.Cells.Item(18, 5).Style = 'Normal'
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '15.99'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  -0.65%  '
$ws.Cells.Item(19, 5).Style = 'Normal'
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  +5.57%  '
$ws.Cells.Item(20, 5).Style = 'Normal'
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '471.90'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '  -0.43%  '
$ws.Cells.Item(21, 5).Style = 'Normal'
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '0.660'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '  +0.98%  '
$ws.Cells.Item(22, 5).Style = 'Normal'
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '79.90'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '  -0.72%  '
$ws.Cells.Item(23, 5).Style = 'Normal'
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '3.833.11'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  -0.53%  '
$ws.Cells.Item(24, 5).Style = 'Normal'
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  +0.54%  '
$ws.Cells.Item(25, 5).Style = 'Normal'
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  -0.07%  '
$ws.Cells.Item(26, 5).Style = 'Normal'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '11.22'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  +1.97%  '
$ws.Cells.Item(27, 5).Style = 'Normal'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '8.93'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  -2.01%  '
$ws.Cells.Item(28, 5).Style = 'Normal'
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  -1.76%  '
$ws.Cells.Item(29, 5).Style = 'Normal'
$ws.Cells.Item(30, 5).NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '  -2.33%  '
$ws.Cells.Item(30, 5).Style = 'Normal'
$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  -1.00%  '
$ws.Cells.Item(31, 5).Style = 'Normal'
$ws.Cells.Item(32, 5).NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '  +0.38%  '
$ws.Cells.Item(32, 5).Style = 'Normal'
$ws.Cells.Item(33, 2).NumberFormat = '@'
$ws.Cells.Item(33, 2).Value = 'NEARProtocol'
$ws.Cells.Item(33, 2).Style = 'Normal'
$ws.Cells.Item(33, 3).NumberFormat = '@'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(33, 3).Style = 'Normal'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '6.55'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '  -0.74%  '
$ws.Cells.Item(33, 5).Style = 'Normal'
$ws.Cells.Item(34, 2).NumberFormat = '@'
$ws.Cells.Item(34, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(34, 2).Style = 'Normal'
$ws.Cells.Item(34, 3).NumberFormat = '@'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(34, 3).Style = 'Normal'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.00'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  -0.01%  '
$ws.Cells.Item(34, 5).Style = 'Normal'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '26.84'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  -0.51%  '
$ws.Cells.Item(35, 5).Style = 'Normal'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '3.682.39'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  -0.35%  '
$ws.Cells.Item(36, 5).Style = 'Normal'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '8.44'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  -0.76%  '
$ws.Cells.Item(37, 5).Style = 'Normal'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '5.91'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  -4.85%  '
$ws.Cells.Item(39, 5).Style = 'Normal'
$ws.Cells.Item(40, 2).NumberFormat = '@'
$ws.Cells.Item(40, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(40, 2).Style = 'Normal'
$ws.Cells.Item(40, 3).NumberFormat = '@'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(40, 3).Style = 'Normal'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '1.00'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  -0.06%  '
$ws.Cells.Item(40, 5).Style = 'Normal'
$ws.Cells.Item(41, 2).NumberFormat = '@'
$ws.Cells.Item(41, 2).Value = 'Monero'
$ws.Cells.Item(41, 2).Style = 'Normal'
$ws.Cells.Item(41, 3).NumberFormat = '@'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(41, 3).Style = 'Normal'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '177.54'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  +6.44%  '
$ws.Cells.Item(41, 5).Style = 'Normal'
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  -1.35%  '
$ws.Cells.Item(42, 5).Style = 'Normal'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.0903'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  -0.24%  '
$ws.Cells.Item(43, 5).Style = 'Normal'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.933'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  -1.33%  '
$ws.Cells.Item(44, 5).Style = 'Normal'
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  +1.77%  '
$ws.Cells.Item(45, 5).Style = 'Normal'
$ws.Cells.Item(46, 2).NumberFormat = '@'
$ws.Cells.Item(46, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(46, 2).Style = 'Normal'
$ws.Cells.Item(46, 3).NumberFormat = '@'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(46, 3).Style = 'Normal'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '29.14'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  +3.25%  '
$ws.Cells.Item(46, 5).Style = 'Normal'
$ws.Cells.Item(47, 2).NumberFormat = '@'
$ws.Cells.Item(47, 2).Value = 'OKB'
$ws.Cells.Item(47, 2).Style = 'Normal'
$ws.Cells.Item(47, 3).NumberFormat = '@'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(47, 3).Style = 'Normal'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '46.69'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  -0.66%  '
$ws.Cells.Item(47, 5).Style = 'Normal'
$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  -2.44%  '
$ws.Cells.Item(48, 5).Style = 'Normal'
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '  -0.03%  '
$ws.Cells.Item(49, 5).Style = 'Normal'
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  -3.68%  '
$ws.Cells.Item(50, 5).Style = 'Normal'
$ws.Cells.Item(51, 5).NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  -5.18%  '
$ws.Cells.Item(51, 5).Style = 'Normal'
